# Updated cryptos list on Thu Mar 30 19:15:57 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $value)
    # Force the cell to stay plain text (matches inline/shared string cells
    # in the source file) instead of letting Excel auto-coerce
    # numeric-looking strings (e.g. "315.77") into a Number cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "27.942.35"
Set-TextCell "E2" "  -1.35%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.778.27"
Set-TextCell "E3" "  -1.39%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.10%  "

# Row 5 - BNB
Set-TextCell "D5" "315.77"
Set-TextCell "E5" "  -0.13%  "

# Row 6 - USDC
Set-TextCell "E6" "  +0.06%  "

# Row 7 - XRP
Set-TextCell "D7" "0.5384"
Set-TextCell "E7" "  -2.74%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.3763"
Set-TextCell "E8" "  -2.41%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.07434"
Set-TextCell "E9" "  -2.33%  "

# Row 10 - OKB
Set-TextCell "D10" "41.63"
Set-TextCell "E10" "  -2.43%  "

# Row 11 - Polygon
Set-TextCell "D11" "1.092"
Set-TextCell "E11" "  -2.47%  "

# Row 12 - BinanceUSD
Set-TextCell "E12" "  +0.07%  "

# Row 13 - Solana
Set-TextCell "D13" "20.41"
Set-TextCell "E13" "  -3.59%  "

# Row 14 - Polkadot
Set-TextCell "D14" "6.067"
Set-TextCell "E14" "  -1.87%  "

# Row 15 / Row 16 - Chainlink and WrappedEther swap positions
Set-TextCell "B15" "WrappedEther"
Set-TextCell "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D15" "1.778.01"
Set-TextCell "E15" "  -1.29%  "

Set-TextCell "B16" "Chainlink"
Set-TextCell "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D16" "7.192"
Set-TextCell "E16" "  -2.29%  "

# Row 17 - Litecoin
Set-TextCell "D17" "88.03"
Set-TextCell "E17" "  -4.60%  "

# Row 18 - ShibaInu
Set-TextCell "E18" "  -1.47%  "

# Row 19 - TRON
Set-TextCell "D19" "0.06430"
Set-TextCell "E19" "  -0.15%  "

# Row 21 - Avalanche
Set-TextCell "E21" "  -0.79%  "

# Row 22 - Uniswap
Set-TextCell "E22" "  -1.99%  "

# Row 23 - WrappedBTC
Set-TextCell "D23" "27.969.10"
Set-TextCell "E23" "  -1.32%  "

# Row 24 - Cosmos
Set-TextCell "D24" "11.13"
Set-TextCell "E24" "  -2.80%  "

# Row 25 - Toncoin
Set-TextCell "D25" "2.082"
Set-TextCell "E25" "  -2.09%  "

# Row 26 - Monero
Set-TextCell "D26" "155.76"
Set-TextCell "E26" "  -1.33%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "20.19"
Set-TextCell "E27" "  -2.25%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-TextCell "D28" "1.977.13"
Set-TextCell "E28" "  -1.71%  "

# Row 29 - LidoDAOToken
Set-TextCell "D29" "2.276"
Set-TextCell "E29" "  -4.99%  "

# Row 30 - BitcoinCash
Set-TextCell "E30" "  -3.46%  "

# Row 31 - ImmutableX
Set-TextCell "E31" "  -1.19%  "

# Row 32 - Stellar
Set-TextCell "D32" "0.1049"
Set-TextCell "E32" "  +2.63%  "

# Row 33 - HuobiToken
Set-TextCell "D33" "3.640"
Set-TextCell "E33" "  -0.73%  "

# Row 34 - Filecoin
Set-TextCell "D34" "5.508"
Set-TextCell "E34" "  -4.00%  "

# Row 35 - Algorand
Set-TextCell "D35" "0.2253"
Set-TextCell "E35" "  -3.59%  "

# Row 36 - Hedera
Set-TextCell "D36" "0.06378"
Set-TextCell "E36" "  +0.83%  "

# Row 37 - VeChain
Set-TextCell "E37" "  -2.54%  "

# Row 38 - InternetComputer(DFINITY)
Set-TextCell "D38" "4.962"
Set-TextCell "E38" "  -1.90%  "

# Row 39 - FraxShare
Set-TextCell "D39" "8.389"
Set-TextCell "E39" "  -5.61%  "

# Row 40 - TheSandbox
Set-TextCell "D40" "0.6118"
Set-TextCell "E40" "  -4.46%  "

# Row 41 - Aptos
Set-TextCell "D41" "11.02"
Set-TextCell "E41" "  -5.28%  "

# Row 42 / Row 43 - TrustWalletToken and WEMIXTOKEN swap positions
Set-TextCell "B42" "WEMIXTOKEN"
Set-TextCell "C42" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D42" "1.428"
Set-TextCell "E42" "  +3.33%  "

Set-TextCell "B43" "TrustWalletToken"
Set-TextCell "C43" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D43" "1.175"
Set-TextCell "E43" "  +1.87%  "

# Row 44 - Frax
Set-TextCell "D44" "0.9997"
Set-TextCell "E44" "  +0.04%  "

# Row 45 - EnergySwap
Set-TextCell "D45" "13.25"
Set-TextCell "E45" "  -1.70%  "

# Row 46 - PancakeSwap
Set-TextCell "D46" "3.654"
Set-TextCell "E46" "  -0.84%  "

# Row 47 - Decentraland
Set-TextCell "D47" "0.5728"
Set-TextCell "E47" "  -4.23%  "

# Row 48 - Quant
Set-TextCell "D48" "126.24"
Set-TextCell "E48" "  +1.51%  "

# Row 49 - EOS
Set-TextCell "D49" "1.182"
Set-TextCell "E49" "  +2.97%  "

# Row 50 - NEARProtocol
Set-TextCell "E50" "  -2.72%  "

# Row 51 - Cronos
Set-TextCell "D51" "0.06785"
Set-TextCell "E51" "  -1.73%  "
